$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.899.75"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.408.87"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "488.91"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Value = "154.84"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  +18.96%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "2.426.31"
$ws.Range("E9").Value = "  -4.12%  "
$ws.Range("D10").Value = "6.27"
$ws.Range("E10").Value = "  +8.68%  "
$ws.Range("D11").Value = "0.0998"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").Value = "2.829.21"
$ws.Range("E14").Value = "  -3.92%  "
$ws.Range("D15").Value = "56.951.08"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("E16").Value = "  -3.79%  "
$ws.Range("E17").Value = "  -3.67%  "
$ws.Range("D18").Value = "2.425.36"
$ws.Range("E18").Value = "  -3.92%  "
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("D20").Value = "324.37"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "5.92"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "57.65"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").Value = "2.511.44"
$ws.Range("E28").Value = "  -3.96%  "
$ws.Range("E29").Value = "  -5.13%  "
$ws.Range("D30").Value = "0.0₃0781"
$ws.Range("E30").Value = "  -5.89%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "150.73"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").Value = "18.56"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").Value = "5.28"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").Value = "3.77"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").Value = "0.831"
$ws.Range("E38").Value = "  -6.35%  "
$ws.Range("E39").Value = "  +8.44%  "
$ws.Range("D40").Value = "34.04"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("D43").Value = "0.994"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "276.62"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "0.597"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("E46").Value = "  -5.97%  "
$ws.Range("D47").Value = "10.22"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("E49").Value = "  -8.12%  "
$ws.Range("D50").Value = "1.891.42"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Value = "17.53"
$ws.Range("E51").Value = "  -3.14%  "
